$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 650, shifting existing rows 650:718 down to 651:719.
$ws.Rows.Item(650).Insert()

# Populate the newly inserted row with the new data record.
$ws.Cells.Item(650, 1).Value = 3
$ws.Cells.Item(650, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(650, 3).Value = "Coquimbo"
$ws.Cells.Item(650, 4).Value = 45194
$ws.Cells.Item(650, 5).Value = 5
$ws.Cells.Item(650, 6).Value = 100112017
$ws.Cells.Item(650, 7).Value = "Apio"
$ws.Cells.Item(650, 8).Value = "Americana (o)"
$ws.Cells.Item(650, 9).Value = "Primera"
$ws.Cells.Item(650, 10).Value = 220
$ws.Cells.Item(650, 11).Value = 8500
$ws.Cells.Item(650, 12).Value = 9000
$ws.Cells.Item(650, 13).Value = 8750
$ws.Cells.Item(650, 14).Value = "`$/docena de matas"
$ws.Cells.Item(650, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(650, 16).Value = 1458
$ws.Cells.Item(650, 17).Value = 6
$ws.Cells.Item(650, 18).Value = "Hortaliza"
